{"js": "// The bibliography paragraph currently holds one run whose <w:t> runs all of\n// \"Bibliografia b\u00e1sica: ... 623p.\" together with no breaks between the\n// citations. We need to turn it into a series of <w:t>/<w:br/> chunks\n// (still inside the same run) so each reference - and the blank line\n// between \"b\u00e1sica\" and \"complementar\" - sits on its own visual line.\nconst LINE_BREAK = \"\\u000b\"; // maps to <w:br/> when written through Range.Text / Paragraph.insertText\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.indexOf(\"Bibliografia b\u00e1sica:\") === 0\n);\nif (!target) {\n  throw new Error(\"Could not find the 'Bibliografia b\u00e1sica:' paragraph\");\n}\n\nconst newText =\n  \"Bibliografia b\u00e1sica:\" +\n  LINE_BREAK +\n  \"PRESS, F.; SIEVER, R.; GROTZINGER, J.; JORDAN, T. H. Para entender a Terra. Porto Alegre: Bookman, 2008. 656p.\" +\n  LINE_BREAK +\n  \"REED, W.; MONROE, J. S. Fundamentos de Geologia. S\u00e3o Paulo: Cengage Learning, 2011. 508p.\" +\n  LINE_BREAK +\n  LINE_BREAK +\n  \"Bibliografia complementar:\" +\n  LINE_BREAK +\n  \"TEIXEIRA, W.; FAIRCHILD, T. R.; DE TOLEDO, M. C. M.; TAIOLI, F. Decifrando a Terra. S\u00e3o Paulo: Companhia Editora Nacional, 2003. 623p.\";\n\ntarget.insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "# The bibliography paragraph currently holds one run whose text runs all of\n# \"Bibliografia basica: ... 623p.\" together with no breaks between the\n# citations. Rewrite it so each reference - and the blank line between\n# \"basica\" and \"complementar\" - sits on its own line, using Word's inline\n# line-break character (vertical tab, Chr(11)) which serializes as <w:br/>\n# while keeping everything inside the same run.\n$d = $word.ActiveDocument\n\n$lineBreak = [char]11\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"Bibliografia b\u00e1sica:\")) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Bibliografia b\u00e1sica:' paragraph\"\n}\n\n$newText = \"Bibliografia b\u00e1sica:\" + $lineBreak +\n    \"PRESS, F.; SIEVER, R.; GROTZINGER, J.; JORDAN, T. H. Para entender a Terra. Porto Alegre: Bookman, 2008. 656p.\" + $lineBreak +\n    \"REED, W.; MONROE, J. S. Fundamentos de Geologia. S\u00e3o Paulo: Cengage Learning, 2011. 508p.\" + $lineBreak + $lineBreak +\n    \"Bibliografia complementar:\" + $lineBreak +\n    \"TEIXEIRA, W.; FAIRCHILD, T. R.; DE TOLEDO, M. C. M.; TAIOLI, F. Decifrando a Terra. S\u00e3o Paulo: Companhia Editora Nacional, 2003. 623p.\"\n\n$target.Range.Text = $newText\n"}
